$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the URL text from E2 to F2 (new destination seller URL), clear E2's value.
$newUrl = "https://www.ozon.ru/seller/ip-yartseva-yu-s-260199/products/?miniapp=seller_260199"

$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = $newUrl

# Increase row 2 height to accommodate the longer wrapped text.
$ws.Rows(2).RowHeight = 165.75
